# Natmi following Dr Hou advice
# Rebuilds the LR-pair data block (rows 2+) as a full 3x3 cross-join of the
# sending/target clusters (ECs, FAPs, sCs) for the fixed Oxt -> Avpr1a pair,
# replacing the previous partial (4-row) table with the new 9-row table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns are, in order: A..T
#   Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
#   Ligand-expressing cells, Ligand detection rate,
#   Ligand average expression value, Ligand total expression value,
#   Ligand derived specificity of average expression value,
#   Ligand derived specificity of total expression value,
#   Receptor-expressing cells, Receptor detection rate,
#   Receptor average expression value, Receptor total expression value,
#   Receptor derived specificity of average expression value,
#   Receptor derived specificity of total expression value,
#   Edge average expression weight, Edge total expression weight,
#   Edge average expression derived specificity, Edge total expression derived specificity
$newRows = @(
  @("ECs", "Oxt", "Avpr1a", "ECs", 3, 1, 0.5822833333333334, 1.74685, 0.5190008657689688, 0.5190008657689686, 2, 0.6666666666666666, 0.3710656666666667, 1.113197, 0.03569378523095436, 0.03569378523095435, 0.2160653532722223, 1.94458817945, 0.01852510543743694, 0.01852510543743693),
  @("ECs", "Oxt", "Avpr1a", "FAPs", 3, 1, 0.5822833333333334, 1.74685, 0.5190008657689688, 0.5190008657689686, 3, 1, 4.128773333333333, 12.38632, 0.3971575973362078, 0.3971575973362078, 2.404115899111111, 21.637043092, 0.2061251368642153, 0.2061251368642153),
  @("ECs", "Oxt", "Avpr1a", "sCs", 3, 1, 0.5822833333333334, 1.74685, 0.5190008657689688, 0.5190008657689686, 3, 1, 5.895967, 17.687901, 0.5671486174328377, 0.5671486174328378, 3.433123317983334, 30.89810986185, 0.2943506234673164, 0.2943506234673164),
  @("FAPs", "Oxt", "Avpr1a", "ECs", 3, 1, 0.4682733333333333, 1.40482, 0.4173814559060953, 0.4173814559060952, 2, 0.6666666666666666, 0.3710656666666667, 1.113197, 0.03569378523095436, 0.03569378523095435, 0.1737601566155555, 1.56384140954, 0.01489792404649521, 0.01489792404649521),
  @("FAPs", "Oxt", "Avpr1a", "FAPs", 3, 1, 0.4682733333333333, 1.40482, 0.4173814559060953, 0.4173814559060952, 3, 1, 4.128773333333333, 12.38632, 0.3971575973362078, 0.3971575973362078, 1.933394451377778, 17.4005500624, 0.1657662162003532, 0.1657662162003531),
  @("FAPs", "Oxt", "Avpr1a", "sCs", 3, 1, 0.4682733333333333, 1.40482, 0.4173814559060953, 0.4173814559060952, 3, 1, 5.895967, 17.687901, 0.5671486174328377, 0.5671486174328378, 2.760924120313333, 24.84831708282, 0.2367173156592469, 0.2367173156592469),
  @("sCs", "Oxt", "Avpr1a", "ECs", 1, 0.3333333333333333, 0.07137466666666667, 0.214124, 0.06361767832493612, 0.0636176783249361, 2, 0.6666666666666666, 0.3710656666666667, 1.113197, 0.03569378523095436, 0.03569378523095435, 0.02648468826977778, 0.238362194428, 0.00227075574702221, 0.002270755747022209),
  @("sCs", "Oxt", "Avpr1a", "FAPs", 1, 0.3333333333333333, 0.07137466666666667, 0.214124, 0.06361767832493612, 0.0636176783249361, 3, 1, 4.128773333333333, 12.38632, 0.3971575973362078, 0.3971575973362078, 0.2946898204088889, 2.65220838368, 0.02526624427163937, 0.02526624427163937),
  @("sCs", "Oxt", "Avpr1a", "sCs", 1, 0.3333333333333333, 0.07137466666666667, 0.214124, 0.06361767832493612, 0.0636176783249361, 3, 1, 5.895967, 17.687901, 0.5671486174328377, 0.5671486174328378, 0.4208226793026666, 3.787404113724, 0.03608067830627452, 0.03608067830627452)
)

$r = 2
foreach ($row in $newRows) {
  $c = 1
  foreach ($val in $row) {
    $ws.Cells.Item($r, $c).Value = $val
    $c = $c + 1
  }
  $r = $r + 1
}

Write-Host "Wrote $($newRows.Count) data rows (rows 2..$($r-1))"
